$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "The Coffeeologist Cafe  70/300 Point Cook Rd  Point Cook VIC 3030"
$ws.Range("C2").Value = "11:00am - 11:40am  8/2/2021"

$ws.Range("B3").Value = "The Coffeeologist Cafe  70/300 Point Cook Rd  Point Cook VIC 3030"
$ws.Range("C3").Value = "11:30am - 12:10pm  0/2/2021"

$ws.Range("B4").Value = "Sunbury Square Shopping Centre  2-28 Evans St  Sunbury VIC 3429"
